# Update "想去人数" (column F) counts on both the "展览" and "全部类型"
# worksheets. These two sheets mirror the same underlying data, so the
# same F-column values need to be bumped on each of them.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F.
$updates = @{
    3  = 3091
    7  = 1709
    9  = 92
    10 = 36
    11 = 4
    12 = 1403
    14 = 538
    16 = 60
    20 = 127
    23 = 3283
    24 = 398
    25 = 162
    26 = 339
    27 = 15
    29 = 161
    30 = 106
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
